# Dagboek - add today's journal entry
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Previous entry ("Opnieuw opnemen van project") is now done
$ws.Range("C10").Value = "COMPLETED"

# Add the new row (today's entry) below the existing last row (row 10)
$today = $excel.Evaluate("TODAY()")

$ws.Range("A11").Value = $today
$ws.Range("B11").Value = "Persistence context configureren"
$ws.Range("C11").Value = "ON TRACK"

# Match the formatting used by the rest of the table (vertically centred text)
$ws.Range("B11").VerticalAlignment = -4108
$ws.Range("C11").VerticalAlignment = -4108

# Fix the authors string spelling (LANGHE -> LANGE)
$ws.Range("D1").Value = "PIETER DELOBELLE, ANTON DANNEELS, MATTHIAS DE LANGE"

# Move the active selection to D1, as in the saved workbook
$ws.Range("D1").Select() | Out-Null
